$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Infant_9m_Toys")
$ws.Range("A13").Value = "'22"
$ws.Range("A14").Value = "'64"
$ws.Range("A15").Value = "'45"
$ws.Range("A16").Value = "'65"
$ws.Range("A17").Value = "'59"
$ws.Range("A18").Value = "'49"
$ws.Range("A19").Value = "'16"
$ws.Range("A20").Value = "'86"
$ws.Range("A21").Value = "'85"
$ws.Range("A22").Value = "'73"
$ws.Range("A23").Value = "'02"
$ws.Range("A24").Value = "'87"
$ws.Range("A25").Value = "'44"
$ws.Range("A26").Value = "'76"

$ws = $wb.Worksheets.Item("Mom_9m_Toys")
$ws.Range("A15").Value = "'80"
$ws.Range("A16").Value = "'65"
$ws.Range("A18").Value = "'68"
$ws.Range("A19").Value = "'75"
$ws.Range("A20").Value = "'63"
$ws.Range("A21").Value = "'12"
$ws.Range("A22").Value = "'42"
$ws.Range("A23").Value = "'78"
$ws.Range("A24").Value = "'73"
$ws.Range("A25").Value = "'87"
$ws.Range("A26").Value = "'43"

$ws = $wb.Worksheets.Item("Infant_9m_NoToys")
$ws.Range("A23").Value = "'37"
$ws.Range("A26").Value = "'21"
$ws.Range("A27").Value = "'24"
$ws.Range("A28").Value = "'49"
$ws.Range("A29").Value = "'60"
$ws.Range("A30").Value = "'16"
$ws.Range("A31").Value = "'15"
$ws.Range("A32").Value = "'57"
$ws.Range("A33").Value = "'73"
$ws.Range("A34").Value = "'02"
$ws.Range("A35").Value = "'50"

$ws = $wb.Worksheets.Item("Mom_9m_NoToys")
$ws.Range("A25").Value = "'65"
$ws.Range("A26").Value = "'59"
$ws.Range("A27").Value = "'51"
$ws.Range("A28").Value = "'75"
$ws.Range("A29").Value = "'42"
$ws.Range("A30").Value = "'78"
$ws.Range("A31").Value = "'73"
$ws.Range("A32").Value = "'77"
$ws.Range("A33").Value = "'87"
$ws.Range("A34").Value = "'88"
$ws.Range("A35").Value = "'43"
